# Update Name of Algo
# Apply numeric value corrections to the imputed result data on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -6.331199999999992
$ws.Range("C3").Value = -11.09679999999999
$ws.Range("D5").Value = -8.776199999999996
$ws.Range("E5").Value = 12.33319999999999
$ws.Range("E9").Value = 12.5314
$ws.Range("E11").Value = 13.2099
$ws.Range("C14").Value = -12.0473
$ws.Range("C21").Value = -13.22880000000001
$ws.Range("E21").Value = 12.93539999999999
$ws.Range("C23").Value = -12.10720000000001
$ws.Range("C25").Value = -11.1592
